$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.71895933333333
$ws.Range("H2").Value = 50.156878
$ws.Range("I2").Value = 0.02912144738161902
$ws.Range("J2").Value = 0.03059269312988411
$ws.Range("M2").Value = 1.378421333333333
$ws.Range("N2").Value = 4.135264
$ws.Range("O2").Value = 0.01656231489052403
$ws.Range("P2").Value = 0.01794267551419991
$ws.Range("Q2").Value = 23.04577021619912
$ws.Range("R2").Value = 207.411931945792
$ws.Range("S2").Value = 0.0004823185816022007
$ws.Range("T2").Value = 0.0005489147659350035
$ws.Range("G3").Value = 16.71895933333333
$ws.Range("H3").Value = 50.156878
$ws.Range("I3").Value = 0.02912144738161902
$ws.Range("J3").Value = 0.03059269312988411
$ws.Range("O3").Value = 0.2170932623988173
$ws.Range("P3").Value = 0.2351865659654651
$ws.Range("Q3").Value = 302.0762178353841
$ws.Range("R3").Value = 2718.685960518456
$ws.Range("S3").Value = 0.006322070017851167
$ws.Range("T3").Value = 0.007194990440852722
$ws.Range("G4").Value = 16.71895933333333
$ws.Range("H4").Value = 50.156878
$ws.Range("I4").Value = 0.02912144738161902
$ws.Range("J4").Value = 0.03059269312988411
$ws.Range("M4").Value = 17.58286933333333
$ws.Range("N4").Value = 52.748608
$ws.Range("O4").Value = 0.2112656061941426
$ws.Range("P4").Value = 0.22887321273073
$ws.Range("Q4").Value = 293.9672773473138
$ws.Range("R4").Value = 2645.705496125824
$ws.Range("S4").Value = 0.006152360234328568
$ws.Range("T4").Value = 0.007001847962721909
$ws.Range("G5").Value = 16.71895933333333
$ws.Range("H5").Value = 50.156878
$ws.Range("I5").Value = 0.02912144738161902
$ws.Range("J5").Value = 0.03059269312988411
$ws.Range("M5").Value = 19.2082395
$ws.Range("N5").Value = 38.416479
$ws.Range("O5").Value = 0.2307951156866419
$ws.Range("P5").Value = 0.1666869194070983
$ws.Range("Q5").Value = 321.1417750654271
$ws.Range("R5").Value = 1926.850650392562
$ws.Range("S5").Value = 0.006721087817403216
$ws.Range("T5").Value = 0.005099401774187084
$ws.Range("G6").Value = 16.71895933333333
$ws.Range("H6").Value = 50.156878
$ws.Range("I6").Value = 0.02912144738161902
$ws.Range("J6").Value = 0.03059269312988411
$ws.Range("M6").Value = 26.988955
$ws.Range("N6").Value = 80.966865
$ws.Range("O6").Value = 0.3242837008298742
$ws.Range("P6").Value = 0.3513106263825066
$ws.Range("Q6").Value = 451.2272410941634
$ws.Range("R6").Value = 4061.04516984747
$ws.Range("S6").Value = 0.009443610730433864
$ws.Range("T6").Value = 0.01074753818618739
$ws.Range("I7").Value = 0.2708539632042961
$ws.Range("J7").Value = 0.2845377865576845
$ws.Range("M7").Value = 1.378421333333333
$ws.Range("N7").Value = 4.135264
$ws.Range("O7").Value = 0.01656231489052403
$ws.Range("P7").Value = 0.01794267551419991
$ws.Range("Q7").Value = 214.3450535392321
$ws.Range("R7").Value = 1929.105481853088
$ws.Range("S7").Value = 0.004485968627935961
$ws.Range("T7").Value = 0.005105369175733208
$ws.Range("I8").Value = 0.2708539632042961
$ws.Range("J8").Value = 0.2845377865576845
$ws.Range("O8").Value = 0.2170932623988173
$ws.Range("P8").Value = 0.2351865659654651
$ws.Range("S8").Value = 0.05880057050566984
$ws.Range("T8").Value = 0.06691946490791631
$ws.Range("I9").Value = 0.2708539632042961
$ws.Range("J9").Value = 0.2845377865576845
$ws.Range("M9").Value = 17.58286933333333
$ws.Range("N9").Value = 52.748608
$ws.Range("O9").Value = 0.2112656061941426
$ws.Range("P9").Value = 0.22887321273073
$ws.Range("Q9").Value = 2734.143021069504
$ws.Range("R9").Value = 24607.28718962553
$ws.Range("S9").Value = 0.0572221267264416
$ws.Range("T9").Value = 0.06512307735274799
$ws.Range("I10").Value = 0.2708539632042961
$ws.Range("J10").Value = 0.2845377865576845
$ws.Range("M10").Value = 19.2082395
$ws.Range("N10").Value = 38.416479
$ws.Range("O10").Value = 0.2307951156866419
$ws.Range("P10").Value = 0.1666869194070983
$ws.Range("Q10").Value = 2986.888714255166
$ws.Range("R10").Value = 17921.332285531
$ws.Range("S10").Value = 0.06251177177192095
$ws.Range("T10").Value = 0.04742872709621491
$ws.Range("I11").Value = 0.2708539632042961
$ws.Range("J11").Value = 0.2845377865576845
$ws.Range("M11").Value = 26.988955
$ws.Range("N11").Value = 80.966865
$ws.Range("O11").Value = 0.3242837008298742
$ws.Range("P11").Value = 0.3513106263825066
$ws.Range("Q11").Value = 4196.793001203496
$ws.Range("R11").Value = 37771.13701083146
$ws.Range("S11").Value = 0.08783352557232769
$ws.Range("T11").Value = 0.09996114802507211
$ws.Range("G12").Value = 194.8548433333333
$ws.Range("H12").Value = 584.56453
$ws.Range("I12").Value = 0.3394024086099587
$ws.Range("J12").Value = 0.3565493705749576
$ws.Range("M12").Value = 1.378421333333333
$ws.Range("N12").Value = 4.135264
$ws.Range("O12").Value = 0.01656231489052403
$ws.Range("P12").Value = 0.01794267551419991
$ws.Range("Q12").Value = 268.5920729539911
$ws.Range("R12").Value = 2417.32865658592
$ws.Range("S12").Value = 0.005621289566000442
$ws.Range("T12").Value = 0.006397449661018682
$ws.Range("G13").Value = 194.8548433333333
$ws.Range("H13").Value = 584.56453
$ws.Range("I13").Value = 0.3394024086099587
$ws.Range("J13").Value = 0.3565493705749576
$ws.Range("O13").Value = 0.2170932623988173
$ws.Range("P13").Value = 0.2351865659654651
$ws.Range("Q13").Value = 3520.614706184841
$ws.Range("R13").Value = 31685.53235566356
$ws.Range("S13").Value = 0.07368197615115236
$ws.Range("T13").Value = 0.08385562206267233
$ws.Range("G14").Value = 194.8548433333333
$ws.Range("H14").Value = 584.56453
$ws.Range("I14").Value = 0.3394024086099587
$ws.Range("J14").Value = 0.3565493705749576
$ws.Range("M14").Value = 17.58286933333333
$ws.Range("N14").Value = 52.748608
$ws.Range("O14").Value = 0.2112656061941426
$ws.Range("P14").Value = 0.22887321273073
$ws.Range("Q14").Value = 3426.107249297137
$ws.Range("R14").Value = 30834.96524367424
$ws.Range("S14").Value = 0.07170405559873501
$ws.Range("T14").Value = 0.08160459994061016
$ws.Range("G15").Value = 194.8548433333333
$ws.Range("H15").Value = 584.56453
$ws.Range("I15").Value = 0.3394024086099587
$ws.Range("J15").Value = 0.3565493705749576
$ws.Range("M15").Value = 19.2082395
$ws.Range("N15").Value = 38.416479
$ws.Range("O15").Value = 0.2307951156866419
$ws.Range("P15").Value = 0.1666869194070983
$ws.Range("Q15").Value = 3742.818498481645
$ws.Range("R15").Value = 22456.91099088987
$ws.Range("S15").Value = 0.07833241815946033
$ws.Range("T15").Value = 0.0594321161976796
$ws.Range("G16").Value = 194.8548433333333
$ws.Range("H16").Value = 584.56453
$ws.Range("I16").Value = 0.3394024086099587
$ws.Range("J16").Value = 0.3565493705749576
$ws.Range("M16").Value = 26.988955
$ws.Range("N16").Value = 80.966865
$ws.Range("O16").Value = 0.3242837008298742
$ws.Range("P16").Value = 0.3513106263825066
$ws.Range("Q16").Value = 5258.928598255383
$ws.Range("R16").Value = 47330.35738429845
$ws.Range("S16").Value = 0.1100626691346106
$ws.Range("T16").Value = 0.1252595827129768
$ws.Range("G17").Value = 82.82950199999999
$ws.Range("H17").Value = 165.659004
$ws.Range("I17").Value = 0.1442742299952585
$ws.Range("J17").Value = 0.1010420758958371
$ws.Range("M17").Value = 1.378421333333333
$ws.Range("N17").Value = 4.135264
$ws.Range("O17").Value = 0.01656231489052403
$ws.Range("P17").Value = 0.01794267551419991
$ws.Range("Q17").Value = 114.173952586176
$ws.Range("R17").Value = 685.0437155170559
$ws.Range("S17").Value = 0.002389515227769359
$ws.Range("T17").Value = 0.001812965181080166
$ws.Range("G18").Value = 82.82950199999999
$ws.Range("H18").Value = 165.659004
$ws.Range("I18").Value = 0.1442742299952585
$ws.Range("J18").Value = 0.1010420758958371
$ws.Range("O18").Value = 0.2170932623988173
$ws.Range("P18").Value = 0.2351865659654651
$ws.Range("Q18").Value = 1496.553833913768
$ws.Range("R18").Value = 8979.323003482608
$ws.Range("S18").Value = 0.03132096326974797
$ws.Range("T18").Value = 0.02376373884796383
$ws.Range("G19").Value = 82.82950199999999
$ws.Range("H19").Value = 165.659004
$ws.Range("I19").Value = 0.1442742299952585
$ws.Range("J19").Value = 0.1010420758958371
$ws.Range("M19").Value = 17.58286933333333
$ws.Range("N19").Value = 52.748608
$ws.Range("O19").Value = 0.2112656061941426
$ws.Range("P19").Value = 0.22887321273073
$ws.Range("Q19").Value = 1456.380310611072
$ws.Range("R19").Value = 8738.281863666431
$ws.Range("S19").Value = 0.03048018265814144
$ws.Range("T19").Value = 0.0231258245312625
$ws.Range("G20").Value = 82.82950199999999
$ws.Range("H20").Value = 165.659004
$ws.Range("I20").Value = 0.1442742299952585
$ws.Range("J20").Value = 0.1010420758958371
$ws.Range("M20").Value = 19.2082395
$ws.Range("N20").Value = 38.416479
$ws.Range("O20").Value = 0.2307951156866419
$ws.Range("P20").Value = 0.1666869194070983
$ws.Range("Q20").Value = 1591.008912081729
$ws.Range("R20").Value = 6364.035648326916
$ws.Range("S20").Value = 0.03329778760235687
$ws.Range("T20").Value = 0.01684239236157532
$ws.Range("G21").Value = 82.82950199999999
$ws.Range("H21").Value = 165.659004
$ws.Range("I21").Value = 0.1442742299952585
$ws.Range("J21").Value = 0.1010420758958371
$ws.Range("M21").Value = 26.988955
$ws.Range("N21").Value = 80.966865
$ws.Range("O21").Value = 0.3242837008298742
$ws.Range("P21").Value = 0.3513106263825066
$ws.Range("Q21").Value = 2235.48170215041
$ws.Range("R21").Value = 13412.89021290246
$ws.Range("S21").Value = 0.04678578123724286
$ws.Range("T21").Value = 0.03549715497395532
$ws.Range("G22").Value = 124.2078576666667
$ws.Range("H22").Value = 372.623573
$ws.Range("I22").Value = 0.2163479508088675
$ws.Range("J22").Value = 0.2272780738416368
$ws.Range("M22").Value = 1.378421333333333
$ws.Range("N22").Value = 4.135264
$ws.Range("O22").Value = 0.01656231489052403
$ws.Range("P22").Value = 0.01794267551419991
$ws.Range("Q22").Value = 171.2107607753636
$ws.Range("R22").Value = 1540.896846978272
$ws.Range("S22").Value = 0.003583222887216068
$ws.Range("T22").Value = 0.004077976730432857
$ws.Range("G23").Value = 124.2078576666667
$ws.Range("H23").Value = 372.623573
$ws.Range("I23").Value = 0.2163479508088675
$ws.Range("J23").Value = 0.2272780738416368
$ws.Range("O23").Value = 0.2170932623988173
$ws.Range("P23").Value = 0.2351865659654651
$ws.Range("Q23").Value = 2244.173164209844
$ws.Range("R23").Value = 20197.5584778886
$ws.Range("S23").Value = 0.04696768245439589
$ws.Range("T23").Value = 0.05345274970605998
$ws.Range("G24").Value = 124.2078576666667
$ws.Range("H24").Value = 372.623573
$ws.Range("I24").Value = 0.2163479508088675
$ws.Range("J24").Value = 0.2272780738416368
$ws.Range("M24").Value = 17.58286933333333
$ws.Range("N24").Value = 52.748608
$ws.Range("O24").Value = 0.2112656061941426
$ws.Range("P24").Value = 0.22887321273073
$ws.Range("Q24").Value = 2183.930531526265
$ws.Range("R24").Value = 19655.37478373638
$ws.Range("S24").Value = 0.04570688097649594
$ws.Range("T24").Value = 0.05201786294338752
$ws.Range("G25").Value = 124.2078576666667
$ws.Range("H25").Value = 372.623573
$ws.Range("I25").Value = 0.2163479508088675
$ws.Range("J25").Value = 0.2272780738416368
$ws.Range("M25").Value = 19.2082395
$ws.Range("N25").Value = 38.416479
$ws.Range("O25").Value = 0.2307951156866419
$ws.Range("P25").Value = 0.1666869194070983
$ws.Range("Q25").Value = 2385.814277843244
$ws.Range("R25").Value = 14314.88566705947
$ws.Range("S25").Value = 0.04993205033550049
$ws.Range("T25").Value = 0.03788428197744147
$ws.Range("G26").Value = 124.2078576666667
$ws.Range("H26").Value = 372.623573
$ws.Range("I26").Value = 0.2163479508088675
$ws.Range("J26").Value = 0.2272780738416368
$ws.Range("M26").Value = 26.988955
$ws.Range("N26").Value = 80.966865
$ws.Range("O26").Value = 0.3242837008298742
$ws.Range("P26").Value = 0.3513106263825066
$ws.Range("Q26").Value = 3352.240281212071
$ws.Range("R26").Value = 30170.16253090864
$ws.Range("S26").Value = 0.07015811415525913
$ws.Range("T26").Value = 0.07984520248431502
